# Weekly fruit/vegetable price data update: shuffle the Fecha/Volumen/Precio values
# across existing rows 2-27 (except 17, 19, 22 which stay the same) per the new data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- source row 5
$ws.Range("D2").Value = 44421
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 8000
$ws.Range("L2").Value = 9000
$ws.Range("M2").Value = 8500
$ws.Range("N2").Value = '$/caja 60 unidades'
$ws.Range("P2").Value = 142
$ws.Range("Q2").Value = 60

# Row 3 <- source row 10
$ws.Range("D3").Value = 44963
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 130
$ws.Range("K3").Value = 4000
$ws.Range("L3").Value = 4500
$ws.Range("M3").Value = 4250
$ws.Range("N3").Value = '$/caja 60 unidades'
$ws.Range("P3").Value = 71
$ws.Range("Q3").Value = 60

# Row 4 <- source row 16
$ws.Range("D4").Value = 44382
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 160
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 7438
$ws.Range("N4").Value = '$/caja 60 unidades'
$ws.Range("P4").Value = 124
$ws.Range("Q4").Value = 60

# Row 5 <- source row 23
$ws.Range("D5").Value = 44648
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 6500
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 6750
$ws.Range("N5").Value = '$/caja 60 unidades'
$ws.Range("P5").Value = 112
$ws.Range("Q5").Value = 60

# Row 6 <- source row 8
$ws.Range("D6").Value = 44785
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 130
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 8000
$ws.Range("M6").Value = 7500
$ws.Range("N6").Value = '$/caja 60 unidades'
$ws.Range("P6").Value = 125
$ws.Range("Q6").Value = 60

# Row 7 <- source row 24
$ws.Range("D7").Value = 44935
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 6000
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 6500
$ws.Range("N7").Value = '$/caja 60 unidades'
$ws.Range("P7").Value = 108
$ws.Range("Q7").Value = 60

# Row 8 <- source row 27
$ws.Range("D8").Value = 44589
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 110
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = 5500
$ws.Range("N8").Value = '$/caja 60 unidades'
$ws.Range("P8").Value = 92
$ws.Range("Q8").Value = 60

# Row 9 <- source row 7
$ws.Range("D9").Value = 45177
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 5000
$ws.Range("L9").Value = 5500
$ws.Range("M9").Value = 5250
$ws.Range("N9").Value = '$/caja 60 unidades'
$ws.Range("P9").Value = 88
$ws.Range("Q9").Value = 60

# Row 10 <- source row 20
$ws.Range("D10").Value = 44657
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 5500
$ws.Range("M10").Value = 5250
$ws.Range("N10").Value = '$/caja 60 unidades'
$ws.Range("P10").Value = 88
$ws.Range("Q10").Value = 60

# Row 11 <- source row 3
$ws.Range("D11").Value = 44281
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 5500
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = 5750
$ws.Range("N11").Value = '$/caja 60 unidades'
$ws.Range("P11").Value = 96
$ws.Range("Q11").Value = 60

# Row 12 <- source row 4
$ws.Range("D12").Value = 44242
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 160
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 5500
$ws.Range("M12").Value = 5250
$ws.Range("N12").Value = '$/caja 60 unidades'
$ws.Range("P12").Value = 88
$ws.Range("Q12").Value = 60

# Row 13 <- source row 18
$ws.Range("D13").Value = 44967
$ws.Range("I13").Value = 'Segunda'
$ws.Range("J13").Value = 50
$ws.Range("K13").Value = 4500
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = 4850
$ws.Range("N13").Value = '$/caja 90 unidades'
$ws.Range("P13").Value = 54
$ws.Range("Q13").Value = 90

# Row 14 <- source row 6
$ws.Range("D14").Value = 44494
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 120
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = 5500
$ws.Range("N14").Value = '$/caja 60 unidades'
$ws.Range("P14").Value = 92
$ws.Range("Q14").Value = 60

# Row 15 <- source row 14
$ws.Range("D15").Value = 44827
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 120
$ws.Range("K15").Value = 6000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 6500
$ws.Range("N15").Value = '$/caja 60 unidades'
$ws.Range("P15").Value = 108
$ws.Range("Q15").Value = 60

# Row 16 <- source row 21
$ws.Range("D16").Value = 44760
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 130
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 7500
$ws.Range("M16").Value = 7250
$ws.Range("N16").Value = '$/caja 60 unidades'
$ws.Range("P16").Value = 121
$ws.Range("Q16").Value = 60

# Row 18 <- source row 12
$ws.Range("D18").Value = 44400
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 9000
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = 9500
$ws.Range("N18").Value = '$/caja 60 unidades'
$ws.Range("P18").Value = 158
$ws.Range("Q18").Value = 60

# Row 20 <- source row 9
$ws.Range("D20").Value = 44627
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = 4000
$ws.Range("L20").Value = 4500
$ws.Range("M20").Value = 4250
$ws.Range("N20").Value = '$/caja 60 unidades'
$ws.Range("P20").Value = 71
$ws.Range("Q20").Value = 60

# Row 21 <- source row 25
$ws.Range("D21").Value = 44764
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 120
$ws.Range("K21").Value = 7000
$ws.Range("L21").Value = 8000
$ws.Range("M21").Value = 7500
$ws.Range("N21").Value = '$/caja 60 unidades'
$ws.Range("P21").Value = 125
$ws.Range("Q21").Value = 60

# Row 23 <- source row 2
$ws.Range("D23").Value = 45079
$ws.Range("I23").Value = 'Primera'
$ws.Range("J23").Value = 130
$ws.Range("K23").Value = 4000
$ws.Range("L23").Value = 5000
$ws.Range("M23").Value = 4462
$ws.Range("N23").Value = '$/caja 60 unidades'
$ws.Range("P23").Value = 74
$ws.Range("Q23").Value = 60

# Row 24 <- source row 15
$ws.Range("D24").Value = 44676
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value = 120
$ws.Range("K24").Value = 4000
$ws.Range("L24").Value = 4500
$ws.Range("M24").Value = 4250
$ws.Range("N24").Value = '$/caja 60 unidades'
$ws.Range("P24").Value = 71
$ws.Range("Q24").Value = 60

# Row 25 <- source row 13
$ws.Range("D25").Value = 44669
$ws.Range("I25").Value = 'Primera'
$ws.Range("J25").Value = 130
$ws.Range("K25").Value = 4500
$ws.Range("L25").Value = 5000
$ws.Range("M25").Value = 4750
$ws.Range("N25").Value = '$/caja 60 unidades'
$ws.Range("P25").Value = 79
$ws.Range("Q25").Value = 60

# Row 26 <- source row 11
$ws.Range("D26").Value = 45243
$ws.Range("I26").Value = 'Primera'
$ws.Range("J26").Value = 120
$ws.Range("K26").Value = 7000
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = 7500
$ws.Range("N26").Value = '$/caja 60 unidades'
$ws.Range("P26").Value = 125
$ws.Range("Q26").Value = 60

# Row 27 <- source row 26
$ws.Range("D27").Value = 45044
$ws.Range("I27").Value = 'Primera'
$ws.Range("J27").Value = 190
$ws.Range("K27").Value = 4000
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = 4526
$ws.Range("N27").Value = '$/caja 60 unidades'
$ws.Range("P27").Value = 75
$ws.Range("Q27").Value = 60

